$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that need to remain literal text even though the new value
# looks numeric (Excel would otherwise auto-convert it to a Number).
$textCells = @( "D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51" )
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = '25.564.79'
$ws.Range("E2").Value = '  +2.43%  '
$ws.Range("D3").Value = '1.672.54'
$ws.Range("E3").Value = '  +2.09%  '
$ws.Range("D4").Value = '0.9990'
$ws.Range("E4").Value = '  +0.26%  '
$ws.Range("D5").Value = '240.08'
$ws.Range("E5").Value = '  +1.45%  '
$ws.Range("D6").Value = '0.9999'
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("D7").Value = '0.4765'
$ws.Range("E7").Value = '  +0.74%  '
$ws.Range("D8").Value = '0.2627'
$ws.Range("E8").Value = '  +2.37%  '
$ws.Range("D9").Value = '0.06176'
$ws.Range("E9").Value = '  +2.66%  '
$ws.Range("D10").Value = '1.672.52'
$ws.Range("E10").Value = '  +2.02%  '
$ws.Range("D11").Value = '0.06984'
$ws.Range("E11").Value = '  -0.92%  '
$ws.Range("D12").Value = '14.89'
$ws.Range("E12").Value = '  +0.40%  '
$ws.Range("D13").Value = '0.5917'
$ws.Range("E13").Value = '  -4.15%  '
$ws.Range("D14").Value = '4.381'
$ws.Range("E14").Value = '  +0.20%  '
$ws.Range("D15").Value = '75.43'
$ws.Range("E15").Value = '  +3.61%  '
$ws.Range("D16").Value = '0.9999'
$ws.Range("E16").Value = '  +0.00%  '
$ws.Range("D17").Value = '0.9996'
$ws.Range("E17").Value = '  +0.19%  '
$ws.Range("D18").Value = '25.563.91'
$ws.Range("E18").Value = '  +2.47%  '
$ws.Range("D19").Value = '0.000006759'
$ws.Range("E19").Value = '  +2.59%  '
$ws.Range("D20").Value = '11.44'
$ws.Range("E20").Value = '  +2.46%  '
$ws.Range("D21").Value = '1.887.96'
$ws.Range("E21").Value = '  +2.25%  '
$ws.Range("D22").Value = '4.457'
$ws.Range("E22").Value = '  +0.97%  '
$ws.Range("D23").Value = '8.801'
$ws.Range("E23").Value = '  +2.27%  '
$ws.Range("D24").Value = '5.282'
$ws.Range("E24").Value = '  -0.11%  '
$ws.Range("D25").Value = '136.79'
$ws.Range("E25").Value = '  +2.77%  '
$ws.Range("D26").Value = '15.05'
$ws.Range("E26").Value = '  +1.64%  '
$ws.Range("E27").Value = '  +2.11%  '
$ws.Range("D28").Value = '1.733'
$ws.Range("E28").Value = '  +4.38%  '
$ws.Range("D29").Value = '104.78'
$ws.Range("E29").Value = '  +2.14%  '
$ws.Range("D30").Value = '3.986'
$ws.Range("E30").Value = '  +6.17%  '
$ws.Range("D31").Value = '0.07877'
$ws.Range("E31").Value = '  +1.90%  '
$ws.Range("D32").Value = '3.638'
$ws.Range("E32").Value = '  +2.15%  '
$ws.Range("D33").Value = '0.9988'
$ws.Range("E33").Value = '  +0.04%  '
$ws.Range("D34").Value = '0.04301'
$ws.Range("E34").Value = '  -0.38%  '
$ws.Range("D35").Value = '2.624'
$ws.Range("E35").Value = '  +0.98%  '
$ws.Range("D36").Value = '0.9608'
$ws.Range("E36").Value = '  +4.10%  '
$ws.Range("D37").Value = '0.6087'
$ws.Range("E37").Value = '  +4.55%  '
$ws.Range("D38").Value = '2.600'
$ws.Range("E38").Value = '  +1.12%  '
$ws.Range("D39").Value = '0.8999'
$ws.Range("E39").Value = '  +8.00%  '
$ws.Range("D40").Value = '0.9999'
$ws.Range("E40").Value = '  +0.20%  '
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = '0.01487'
$ws.Range("E41").Value = '  -4.36%  '
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").Value = '1.862'
$ws.Range("E42").Value = '  +3.30%  '
$ws.Range("D43").Value = '96.36'
$ws.Range("D44").Value = '0.3766'
$ws.Range("E44").Value = '  +1.21%  '
$ws.Range("D45").Value = '4.913'
$ws.Range("E45").Value = '  +3.58%  '
$ws.Range("D46").Value = '0.1123'
$ws.Range("E46").Value = '  +1.61%  '
$ws.Range("D47").Value = '6.235'
$ws.Range("E47").Value = '  +2.39%  '
$ws.Range("D48").Value = '0.05270'
$ws.Range("E48").Value = '  +1.01%  '
$ws.Range("D49").Value = '29.95'
$ws.Range("E49").Value = '  +1.15%  '
$ws.Range("D50").Value = '7.427'
$ws.Range("E50").Value = '  +3.75%  '
$ws.Range("B51").Value = 'TrueUSD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd'
$ws.Range("D51").Value = '1.001'
$ws.Range("E51").Value = '  +0.34%  '
